$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.5
$ws.Range("I2").Value = 2.54
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.35
$ws.Range("N2").Value = 4.9
$ws.Range("P2").Value = 2.28
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 1.51
$ws.Range("S2").Value = 2.92
$ws.Range("T2").Value = 1.61
$ws.Range("U2").Value = 2.54
$ws.Range("V2").Value = 1.65
$ws.Range("X2").Value = 19
$ws.Range("AB2").Value = 15.5
$ws.Range("AG2").Value = 13

# Row 3
$ws.Range("G3").Value = 26
$ws.Range("K3").Value = 10.5
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 1.02
$ws.Range("S3").Value = 2.02
$ws.Range("Y3").Value = 150
$ws.Range("Z3").Value = 980
$ws.Range("AA3").Value = 1000

# Row 4
$ws.Range("F4").Value = 13.5
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 1.27
$ws.Range("I4").Value = 1.31
$ws.Range("J4").Value = 5.1
$ws.Range("K4").Value = 6.6
$ws.Range("L4").Value = 1.32
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 3.5
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 1.87
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.33
$ws.Range("S4").Value = 3.1
$ws.Range("T4").Value = 2.74
$ws.Range("U4").Value = 1.48
$ws.Range("V4").Value = 4.1
$ws.Range("W4").Value = 1.05
$ws.Range("Y4").Value = 7.6
$ws.Range("Z4").Value = 7.2
$ws.Range("AA4").Value = 17.5

# Row 5
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 2.62
$ws.Range("I5").Value = 2.86
$ws.Range("K5").Value = 3.7
$ws.Range("L5").Value = 1.33
$ws.Range("N5").Value = 3.55
$ws.Range("P5").Value = 1.87
$ws.Range("Q5").Value = 2.02
$ws.Range("R5").Value = 1.39
$ws.Range("T5").Value = 1.79
$ws.Range("U5").Value = 2.12
$ws.Range("V5").Value = 1.53
$ws.Range("W5").Value = 1.5
$ws.Range("AH5").Value = 18.5

# Row 6
$ws.Range("G6").Value = 5.2
$ws.Range("H6").Value = 1.86
$ws.Range("K6").Value = 4.8
$ws.Range("V6").Value = 1.84
$ws.Range("W6").Value = 1.23
$ws.Range("AA6").Value = 900

# Row 7
$ws.Range("H7").Value = 2.26
$ws.Range("J7").Value = 3.6
$ws.Range("L7").Value = 1.39
$ws.Range("P7").Value = 2.16
$ws.Range("Q7").Value = 1.84
$ws.Range("U7").Value = 2.36
$ws.Range("AC7").Value = 8
$ws.Range("AE7").Value = 21
$ws.Range("AG7").Value = 14
$ws.Range("AL7").Value = 42

# Row 8
$ws.Range("F8").Value = 3.3
$ws.Range("H8").Value = 2.28
$ws.Range("L8").Value = 1.4
$ws.Range("Y8").Value = 8.800000000000001
$ws.Range("AA8").Value = 900

# Row 9
$ws.Range("G9").Value = 3.2
$ws.Range("H9").Value = 2.3
$ws.Range("W9").Value = 1.45

# Row 11
$ws.Range("L11").Value = 1.31
$ws.Range("Q11").Value = 1.91

# Row 12
$ws.Range("G12").Value = 3.55
$ws.Range("H12").Value = 2.48
$ws.Range("N12").Value = 2.9
$ws.Range("T12").Value = 1.94
$ws.Range("Z12").Value = 980
$ws.Range("AG12").Value = 1000
$ws.Range("AH12").Value = 1000

# Row 13
$ws.Range("F13").Value = 1.61
$ws.Range("V13").Value = 1.18
$ws.Range("Y13").Value = 970

# Row 14
$ws.Range("K14").Value = 3.7
$ws.Range("L14").Value = 1.32
$ws.Range("R14").Value = 1.29
$ws.Range("AJ14").Value = 900
$ws.Range("AN14").Value = 44
$ws.Range("AO14").Value = 25

# Row 15
$ws.Range("F15").Value = 1.77
$ws.Range("G15").Value = 1.86
$ws.Range("I15").Value = 5.4
$ws.Range("J15").Value = 3.8
$ws.Range("L15").Value = 1.3
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 4.2
$ws.Range("P15").Value = 2.08
$ws.Range("Q15").Value = 1.79
$ws.Range("R15").Value = 1.44
$ws.Range("S15").Value = 2.96
$ws.Range("T15").Value = 1.75
$ws.Range("U15").Value = 2.16
$ws.Range("W15").Value = 2.16
$ws.Range("X15").Value = 21

# Row 16
$ws.Range("K16").Value = 3.8
$ws.Range("L16").Value = 1.39
$ws.Range("P16").Value = 1.93
$ws.Range("U16").Value = 2.1
$ws.Range("Y16").Value = 14
$ws.Range("AH16").Value = 28

# Row 17
$ws.Range("I17").Value = 3.15
$ws.Range("K17").Value = 3.55
$ws.Range("L17").Value = 1.37
$ws.Range("P17").Value = 1.84

# Row 18
$ws.Range("F18").Value = 1.84
$ws.Range("G18").Value = 1.96
$ws.Range("H18").Value = 4.2
$ws.Range("I18").Value = 5
$ws.Range("Q18").Value = 1.7
$ws.Range("V18").Value = 1.26
$ws.Range("W18").Value = 2.04
$ws.Range("Y18").Value = 1000
$ws.Range("AD18").Value = 1000
$ws.Range("AF18").Value = 500
$ws.Range("AI18").Value = 1000

# Row 19
$ws.Range("F19").Value = 2.16
$ws.Range("G19").Value = 2.18
$ws.Range("I19").Value = 3.75
$ws.Range("O19").Value = 1.26
$ws.Range("Q19").Value = 1.82
$ws.Range("T19").Value = 1.68
$ws.Range("V19").Value = 1.36
$ws.Range("W19").Value = 1.84
$ws.Range("X19").Value = 16.5
$ws.Range("Y19").Value = 15.5
$ws.Range("Z19").Value = 26
$ws.Range("AJ19").Value = 25
$ws.Range("AL19").Value = 30
$ws.Range("AM19").Value = 65

# Row 20
$ws.Range("M20").Value = 1.04
$ws.Range("P20").Value = 2.78
$ws.Range("Q20").Value = 1.55
$ws.Range("U20").Value = 2.34
$ws.Range("AI20").Value = 25
$ws.Range("AL20").Value = 65
$ws.Range("AN20").Value = 70
$ws.Range("AO20").Value = 5.4

# Row 21
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 4.1
$ws.Range("H21").Value = 1.95
$ws.Range("I21").Value = 1.96
$ws.Range("N21").Value = 5.4
$ws.Range("O21").Value = 1.21
$ws.Range("Q21").Value = 1.68
$ws.Range("T21").Value = 1.64
$ws.Range("V21").Value = 2.04
$ws.Range("W21").Value = 1.32
$ws.Range("Z21").Value = 13.5
$ws.Range("AB21").Value = 18.5
$ws.Range("AG21").Value = 16
$ws.Range("AH21").Value = 15.5
$ws.Range("AN21").Value = 32

# Row 22
$ws.Range("F22").Value = 2.52
$ws.Range("G22").Value = 2.56
$ws.Range("J22").Value = 3.5
$ws.Range("N22").Value = 4
$ws.Range("R22").Value = 1.4
$ws.Range("T22").Value = 1.08
$ws.Range("V22").Value = 1.47
$ws.Range("AK22").Value = 27
$ws.Range("AL22").Value = 36

# Row 23
$ws.Range("G23").Value = 2.08
$ws.Range("K23").Value = 3.5
$ws.Range("P23").Value = 1.57
$ws.Range("T23").Value = 2.12
$ws.Range("W23").Value = 1.92
$ws.Range("AD23").Value = 38

# Row 24
$ws.Range("F24").Value = 3.45
$ws.Range("K24").Value = 3.55
$ws.Range("L24").Value = 1.01

# Row 25
$ws.Range("F25").Value = 2.2
$ws.Range("G25").Value = 2.4
$ws.Range("H25").Value = 3.65
$ws.Range("K25").Value = 3.35
$ws.Range("N25").Value = 2.64
$ws.Range("P25").Value = 1.54
$ws.Range("S25").Value = 5.1
$ws.Range("T25").Value = 2.06
$ws.Range("U25").Value = 1.74
$ws.Range("V25").Value = 1.32
$ws.Range("W25").Value = 1.71
